# Auto-generated Excel COM-interop script
# Applies numeric updates to Leve profit sheets per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2214.2856
$ws.Range("I40").Value = 2083.3333
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2083.3333
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1908.3333
$ws.Range("N40").Value = -3350

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 43125
$ws.Range("J105").Value = 43125
$ws.Range("L105").Value = 43125
$ws.Range("N105").Value = -50113

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 250007860
$ws.Range("I111").Value = 500015000
$ws.Range("J111").Value = 700
$ws.Range("K111").Value = 1500045000
$ws.Range("L111").Value = 2100
$ws.Range("M111").Value = -1500041933
$ws.Range("N111").Value = -8234

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1973.1875
$ws.Range("I112").Value = 625
$ws.Range("J112").Value = 2422.5833
$ws.Range("K112").Value = 1875
$ws.Range("L112").Value = 7267.749899999999
$ws.Range("M112").Value = -767
$ws.Range("N112").Value = -9483.749899999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 9091828
$ws.Range("J115").Value = 1574.6666
$ws.Range("L115").Value = 4723.9998
$ws.Range("N115").Value = -7857.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 695
$ws.Range("I118").Value = 690
$ws.Range("K118").Value = 2070
$ws.Range("M118").Value = -413

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2015.8334
$ws.Range("J125").Value = 2272.6
$ws.Range("L125").Value = 20453.4
$ws.Range("N125").Value = -25373.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2563.6667
$ws.Range("I132").Value = 2841.0908
$ws.Range("J132").Value = 1800.75
$ws.Range("K132").Value = 8523.2724
$ws.Range("L132").Value = 5402.25
$ws.Range("M132").Value = -5993.2724
$ws.Range("N132").Value = -10462.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1529.9032
$ws.Range("I45").Value = 1001
$ws.Range("J45").Value = 1781.762
$ws.Range("K45").Value = 1001
$ws.Range("L45").Value = 1781.762
$ws.Range("M45").Value = -624
$ws.Range("N45").Value = -2535.762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1107.9584
$ws.Range("I61").Value = 1107.9584
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1107.9584
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -895.9584
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1528.7778
$ws.Range("I110").Value = 703.6667
$ws.Range("J110").Value = 1941.3334
$ws.Range("K110").Value = 703.6667
$ws.Range("L110").Value = 1941.3334
$ws.Range("M110").Value = 1341.3333
$ws.Range("N110").Value = -6031.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1692.25
$ws.Range("J122").Value = 1585.4
$ws.Range("L122").Value = 4756.200000000001
$ws.Range("N122").Value = -9656.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1107.9584
$ws.Range("I136").Value = 1107.9584
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3323.8752
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -773.8751999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 112.71429
$ws.Range("I22").Value = 77.8
$ws.Range("K22").Value = 77.8
$ws.Range("M22").Value = 95.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1731.16
$ws.Range("I86").Value = 1565.7142
$ws.Range("K86").Value = 1565.7142
$ws.Range("M86").Value = -442.7141999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1731.16
$ws.Range("I89").Value = 1565.7142
$ws.Range("K89").Value = 7828.571
$ws.Range("M89").Value = -2212.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2816.86
$ws.Range("I105").Value = 1633.3334
$ws.Range("J105").Value = 2853.4639
$ws.Range("K105").Value = 1633.3334
$ws.Range("L105").Value = 2853.4639
$ws.Range("M105").Value = 113.6666
$ws.Range("N105").Value = -6347.463900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4249.1
$ws.Range("I107").Value = 4382.2
$ws.Range("J107").Value = 4116
$ws.Range("K107").Value = 4382.2
$ws.Range("L107").Value = 4116
$ws.Range("M107").Value = -2462.2
$ws.Range("N107").Value = -7956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 317.07693
$ws.Range("I22").Value = 259.45456
$ws.Range("J22").Value = 634
$ws.Range("K22").Value = 259.45456
$ws.Range("L22").Value = 634
$ws.Range("M22").Value = 90.54543999999999
$ws.Range("N22").Value = -1334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4653.5
$ws.Range("I31").Value = 4440
$ws.Range("J31").Value = 5009.3335
$ws.Range("K31").Value = 4440
$ws.Range("L31").Value = 5009.3335
$ws.Range("M31").Value = -4145
$ws.Range("N31").Value = -5599.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4653.5
$ws.Range("I34").Value = 4440
$ws.Range("J34").Value = 5009.3335
$ws.Range("K34").Value = 4440
$ws.Range("L34").Value = 5009.3335
$ws.Range("M34").Value = -4238
$ws.Range("N34").Value = -5413.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4086.2942
$ws.Range("I132").Value = 2682
$ws.Range("J132").Value = 5334.5557
$ws.Range("K132").Value = 8046
$ws.Range("L132").Value = 16003.6671
$ws.Range("M132").Value = -5516
$ws.Range("N132").Value = -21063.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 54858.2
$ws.Range("J140").Value = 54858.2
$ws.Range("L140").Value = 54858.2
$ws.Range("N140").Value = -65218.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 6934.206
$ws.Range("I122").Value = 8437.519
$ws.Range("J122").Value = 1135.7142
$ws.Range("K122").Value = 75937.671
$ws.Range("L122").Value = 10221.4278
$ws.Range("M122").Value = -73487.671
$ws.Range("N122").Value = -15121.4278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3368169.8
$ws.Range("I132").Value = 999
$ws.Range("J132").Value = 5051755
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 45465795
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -45470855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4349659.5
$ws.Range("I122").Value = 5883942.5
$ws.Range("J122").Value = 2524.8333
$ws.Range("K122").Value = 17651827.5
$ws.Range("L122").Value = 7574.499899999999
$ws.Range("M122").Value = -17649377.5
$ws.Range("N122").Value = -12474.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4682.4
$ws.Range("I132").Value = 4335.3335
$ws.Range("K132").Value = 13006.0005
$ws.Range("M132").Value = -10476.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3444.4666
$ws.Range("J7").Value = 3100
$ws.Range("L7").Value = 3100
$ws.Range("N7").Value = -3324

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 690.63635
$ws.Range("I55").Value = 588.55554
$ws.Range("J55").Value = 1150
$ws.Range("K55").Value = 588.55554
$ws.Range("L55").Value = 1150
$ws.Range("M55").Value = -415.55554
$ws.Range("N55").Value = -1496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3444.4666
$ws.Range("J126").Value = 3100
$ws.Range("L126").Value = 9300
$ws.Range("N126").Value = -14240

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4749.75
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5199.6
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 15598.8
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -20658.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 460
$ws.Range("I107").Value = 460
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1380
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 540
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1523.7894
$ws.Range("I132").Value = 1187.129
$ws.Range("J132").Value = 3014.7144
$ws.Range("K132").Value = 3561.387
$ws.Range("L132").Value = 9044.143199999999
$ws.Range("M132").Value = -1031.387
$ws.Range("N132").Value = -14104.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7094.8335
$ws.Range("I136").Value = 7043.875
$ws.Range("J136").Value = 7502.5
$ws.Range("K136").Value = 21131.625
$ws.Range("L136").Value = 22507.5
$ws.Range("M136").Value = -18581.625
$ws.Range("N136").Value = -27607.5
